# Weekly update: a new week of price data (2022-05-24, serial 44705) is
# inserted at the top of the existing "Primera"/"Segunda" record pairs
# (rows 22-23), pushing every following pair down by one (old rows 22-47
# become rows 24-49). The new pair's non-varying columns mirror the pair
# that is being pushed down (same market / region / product metadata).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before row 22, shifting rows 22:47 down
# to 24:49 (this also grows the sheet's used range from R47 to R49, and
# carries over per-column formatting such as the date style on column D).
$ws.Rows.Item(22).Resize(2).Insert()

# New row 22: "Primera" quality for the new week (2022-05-24 -> 44705).
$ws.Cells.Item(22, 1).Value = 11
$ws.Cells.Item(22, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(22, 3).Value = "Bíobío"
$ws.Cells.Item(22, 4).Value = 44705
$ws.Cells.Item(22, 5).Value = 8
$ws.Cells.Item(22, 6).Value = 100112037
$ws.Cells.Item(22, 7).Value = "Cebollín"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 200
$ws.Cells.Item(22, 11).Value = 600
$ws.Cells.Item(22, 12).Value = 700
$ws.Cells.Item(22, 13).Value = 650
$ws.Cells.Item(22, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(22, 15).Value = "Región de Ñuble"
$ws.Cells.Item(22, 16).Value = 108
$ws.Cells.Item(22, 17).Value = 6
$ws.Cells.Item(22, 18).Value = "Hortaliza"

# New row 23: "Segunda" quality for the same new week.
$ws.Cells.Item(23, 1).Value = 11
$ws.Cells.Item(23, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(23, 3).Value = "Bíobío"
$ws.Cells.Item(23, 4).Value = 44705
$ws.Cells.Item(23, 5).Value = 8
$ws.Cells.Item(23, 6).Value = 100112037
$ws.Cells.Item(23, 7).Value = "Cebollín"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Segunda"
$ws.Cells.Item(23, 10).Value = 100
$ws.Cells.Item(23, 11).Value = 500
$ws.Cells.Item(23, 12).Value = 500
$ws.Cells.Item(23, 13).Value = 500
$ws.Cells.Item(23, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(23, 15).Value = "Región de Ñuble"
$ws.Cells.Item(23, 16).Value = 83
$ws.Cells.Item(23, 17).Value = 6
$ws.Cells.Item(23, 18).Value = "Hortaliza"
